$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $scaled = [decimal]$cell.Value2 * [decimal]10000
    $cell.Value2 = ($scaled.ToString() -as [double])
}
